$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "19.954.31"
$ws.Range("E2").Value = "  -8.32%  "
Set-TextValue $ws.Range("D3") "1.418.16"
$ws.Range("E3").Value = "  -7.91%  "
Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  +0.24%  "
Set-TextValue $ws.Range("D5") "1.003"
$ws.Range("E5").Value = "  +0.21%  "
Set-TextValue $ws.Range("D6") "273.39"
$ws.Range("E6").Value = "  -5.90%  "
Set-TextValue $ws.Range("D7") "0.3725"
$ws.Range("E7").Value = "  -4.04%  "
Set-TextValue $ws.Range("D8") "0.3076"
$ws.Range("E8").Value = "  -3.55%  "
Set-TextValue $ws.Range("D9") "39.67"
$ws.Range("E9").Value = "  -8.05%  "
Set-TextValue $ws.Range("D10") "1.010"
$ws.Range("E10").Value = "  -4.76%  "
Set-TextValue $ws.Range("D11") "0.06597"
$ws.Range("E11").Value = "  -8.46%  "
Set-TextValue $ws.Range("D12") "1.003"
$ws.Range("E12").Value = "  +0.27%  "
Set-TextValue $ws.Range("D13") "5.406"
$ws.Range("E13").Value = "  -4.24%  "
Set-TextValue $ws.Range("D14") "17.05"
$ws.Range("E14").Value = "  -8.39%  "
Set-TextValue $ws.Range("D15") "6.155"
Set-TextValue $ws.Range("D16") "1.424.32"
Set-TextValue $ws.Range("D17") "0.00001005"
$ws.Range("E17").Value = "  -9.84%  "
Set-TextValue $ws.Range("D18") "0.05833"
$ws.Range("E18").Value = "  -11.46%  "
Set-TextValue $ws.Range("D19") "74.49"
$ws.Range("E19").Value = "  -10.73%  "
Set-TextValue $ws.Range("D20") "1.003"
$ws.Range("E20").Value = "  +0.24%  "
Set-TextValue $ws.Range("D21") "5.643"
$ws.Range("E21").Value = "  -8.31%  "
Set-TextValue $ws.Range("D22") "14.44"
$ws.Range("E22").Value = "  -6.31%  "
Set-TextValue $ws.Range("D23") "10.90"
$ws.Range("E23").Value = "  -0.52%  "
Set-TextValue $ws.Range("D24") "2.327"
$ws.Range("E24").Value = "  -2.40%  "
Set-TextValue $ws.Range("D25") "19.958.94"
$ws.Range("E25").Value = "  -8.32%  "
Set-TextValue $ws.Range("D26") "2.291"
$ws.Range("E26").Value = "  -4.51%  "
Set-TextValue $ws.Range("D27") "138.81"
$ws.Range("E27").Value = "  -5.33%  "
Set-TextValue $ws.Range("D28") "16.87"
$ws.Range("E28").Value = "  -8.25%  "
Set-TextValue $ws.Range("D29") "1.582.61"
$ws.Range("E29").Value = "  -7.94%  "
Set-TextValue $ws.Range("D30") "108.92"
$ws.Range("E30").Value = "  -7.32%  "
Set-TextValue $ws.Range("D31") "3.813"
$ws.Range("E31").Value = "  -21.28%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "5.426"
$ws.Range("E32").Value = "  -8.14%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D33") "0.8858"
$ws.Range("E33").Value = "  -9.02%  "
Set-TextValue $ws.Range("D34") "0.07741"
$ws.Range("E34").Value = "  -5.73%  "
Set-TextValue $ws.Range("D35") "8.398"
$ws.Range("E35").Value = "  -6.14%  "
Set-TextValue $ws.Range("D36") "11.26"
$ws.Range("E36").Value = "  +5.04%  "
Set-TextValue $ws.Range("D37") "1.002"
$ws.Range("E37").Value = "  +0.21%  "
Set-TextValue $ws.Range("D38") "4.771"
$ws.Range("E38").Value = "  -7.60%  "
Set-TextValue $ws.Range("D39") "0.05675"
$ws.Range("E39").Value = "  -6.77%  "
Set-TextValue $ws.Range("D40") "0.1909"
$ws.Range("E40").Value = "  -6.43%  "
Set-TextValue $ws.Range("D41") "0.02024"
$ws.Range("E41").Value = "  -8.33%  "
Set-TextValue $ws.Range("D42") "1.091"
$ws.Range("E42").Value = "  -8.67%  "
Set-TextValue $ws.Range("D43") "1.264"
$ws.Range("E43").Value = "  -14.80%  "
Set-TextValue $ws.Range("D44") "0.5315"
$ws.Range("E44").Value = "  -7.88%  "
Set-TextValue $ws.Range("D45") "3.530"
$ws.Range("E45").Value = "  -5.79%  "
$ws.Range("E46").Value = "  -6.33%  "
Set-TextValue $ws.Range("D47") "0.5127"
$ws.Range("E47").Value = "  -7.34%  "
Set-TextValue $ws.Range("D48") "1.807"
$ws.Range("E48").Value = "  -3.49%  "
Set-TextValue $ws.Range("D49") "109.46"
$ws.Range("E49").Value = "  -7.49%  "
Set-TextValue $ws.Range("D50") "1.047"
$ws.Range("E50").Value = "  -8.61%  "
$ws.Range("E51").Value = "  +0.31%  "
